# Insert a new data row at row 163 (shifting the existing rows 163..276 down to
# 164..277), and populate the new row with its own data, as described by the
# diff: dimension grows from A1:R276 to A1:R277.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 163; everything below (old 163..276)
# shifts down to 164..277.
$ws.Rows("163:163").Insert()

# Populate the newly inserted row 163 with the new record's values.
$ws.Cells.Item(163, 1).Value  = 3
$ws.Cells.Item(163, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(163, 3).Value  = "Coquimbo"
$ws.Cells.Item(163, 4).Value  = 44673
$ws.Cells.Item(163, 5).Value  = 5
$ws.Cells.Item(163, 6).Value  = 100112001
$ws.Cells.Item(163, 7).Value  = "Berenjena"
$ws.Cells.Item(163, 8).Value  = "Sin especificar"
$ws.Cells.Item(163, 9).Value  = "Primera"
$ws.Cells.Item(163, 10).Value = 125
$ws.Cells.Item(163, 11).Value = 7000
$ws.Cells.Item(163, 12).Value = 7500
$ws.Cells.Item(163, 13).Value = 7240
$ws.Cells.Item(163, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(163, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(163, 16).Value = 121
$ws.Cells.Item(163, 17).Value = 60
$ws.Cells.Item(163, 18).Value = "Hortaliza"

# Make sure the date cell keeps the date/time number format used by the rest
# of column D.
$ws.Cells.Item(163, 4).NumberFormat = $ws.Cells.Item(164, 4).NumberFormat
